$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("optimization_parameters")

# Row 1: remove C1:F1 (previously held "value")
$ws.Range("C1:F1").ClearContents()

# Row 8: rename "Model" label to "production_function" (keep B8 = Sigmoid)
$ws.Cells.Item(8,1).Value = "production_function"

# Insert a new row 9 for "L_curve" = 0, pushing old rows 9+ down
$ws.Rows.Item(9).Insert()
$ws.Cells.Item(9,1).Value = "L_curve"
$ws.Cells.Item(9,2).Value = 0
$ws.Cells.Item(9,2).NumberFormat = "0.00E+00"

# Remove the "Deletion" row entirely (now shifted down to row 17)
$ws.Rows.Item(17).Delete()

Write-Host ("A8=" + $ws.Cells.Item(8,1).Value())
Write-Host ("B8=" + $ws.Cells.Item(8,2).Value())
Write-Host ("A9=" + $ws.Cells.Item(9,1).Value())
Write-Host ("B9=" + $ws.Cells.Item(9,2).Value())
Write-Host ("A10=" + $ws.Cells.Item(10,1).Value())
Write-Host ("A15=" + $ws.Cells.Item(15,1).Value())
Write-Host ("A16=" + $ws.Cells.Item(16,1).Value())
Write-Host ("A17=" + $ws.Cells.Item(17,1).Value())
